# Apply "last minute bug fix and testing" changes to results table:
#  - insert KB sub-folder paths ("Sammys/KB/", "Queens/KB/", "Map/KB/")
#    in front of the .cnf filenames used in each "python DPLL.py ..." cell
#  - resize the results table's columns

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Get-CellRange($row, $col) {
    $c = $t.Cell($row, $col)
    # Re-anchoring the cell's start/end on $d.Range(...) (instead of using
    # $c.Range / $c.Range.Duplicate directly) is what makes
    # Find/Collapse/InsertBefore/InsertAfter behave correctly here.
    return $d.Range($c.Range.Start, $c.Range.End)
}

# ---- Sammys A (row 2) ---------------------------------------------------
$rng = Get-CellRange 2 2
$rng.Find.Execute("sammy.cnf") | Out-Null
$rng.Collapse(1)
$rng.InsertBefore("Sammys/KB/")

# ---- Sammys B (row 3) ---------------------------------------------------
$rng = Get-CellRange 3 2
$rng.Find.Execute("sammy.cnf") | Out-Null
$rng.Collapse(1)
$rng.InsertBefore("Sammys/KB/")

# ---- 3 Queens (row 4) -----------------------------------------------------
$rng = Get-CellRange 4 2
$rng.Find.Execute("3queens.cnf") | Out-Null
$rng.Collapse(1)
$rng.InsertBefore("Queens/KB/")

# ---- 4 Queens (row 5) -----------------------------------------------------
$rng = Get-CellRange 5 2
$rng.Find.Execute("4queens.cnf") | Out-Null
$rng.Collapse(1)
$rng.InsertBefore("Queens/KB/")

# ---- 5 Queens (row 6) -----------------------------------------------------
$rng = Get-CellRange 6 2
$rng.Find.Execute("5queens.cnf") | Out-Null
$rng.Collapse(1)
$rng.InsertBefore("Queens/KB/")

# ---- 6 Queens (row 7) -----------------------------------------------------
$rng = Get-CellRange 7 2
$rng.Find.Execute("6queens.cnf") | Out-Null
$rng.Collapse(1)
$rng.InsertBefore("Queens/KB/")

# ---- Map Color (row 8) -----------------------------------------------------
$rng = Get-CellRange 8 2
$rng.Find.Execute("mapcolor.cnf") | Out-Null
$rng.Collapse(1)
$rng.InsertBefore("Map/KB/")

# ---- Map Color 2 (row 9) ---------------------------------------------------
$rng = Get-CellRange 9 2
$rng.Find.Execute("mapcolor.cnf") | Out-Null
$rng.Collapse(1)
$rng.InsertBefore("Map/KB/")

# ---- Map Color Unsatisfiable (row 10) --------------------------------------
$rng = Get-CellRange 10 2
$rng.Find.Execute("mapcolor.cnf") | Out-Null
$rng.Collapse(1)
$rng.InsertBefore("Map/KB/")

# ---- Resize table columns --------------------------------------------------
# 1870/1870/1870/1870/1870 twips -> 1794/2399/1793/1682/1682 twips
# Column.Width is expressed in points (1 twip = 1/20 pt).
$t.Columns.Item(1).Width = 1794 / 20.0
$t.Columns.Item(2).Width = 2399 / 20.0
$t.Columns.Item(3).Width = 1793 / 20.0
$t.Columns.Item(4).Width = 1682 / 20.0
$t.Columns.Item(5).Width = 1682 / 20.0

Write-Host "Edit complete"
